$wb = $excel.ActiveWorkbook

# Updated "想去人数" (interested-count) values for the 展览 and 全部类型 sheets.
$updates = @{
    "F2"  = 6656
    "F3"  = 189
    "F5"  = 49
    "F6"  = 2046
    "F7"  = 1559
    "F9"  = 1015
    "F10" = 442
    "F11" = 17
    "F12" = 5639
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
